# Make it possible to select different features to be executed.
# Currently only "Check Best Practices" is implemented, so every check
# filename gets a "CheckBestPractices\" folder prefix, and the "Project"
# sheet's Action/Fix column is aligned with the "Workflow" sheet's wording.

$wb = $excel.ActiveWorkbook

# --- Workflow sheet -------------------------------------------------
$wsWorkflow = $wb.Worksheets.Item("Workflow")

for ($row = 2; $row -le 12; $row++) {
    $cell = $wsWorkflow.Cells.Item($row, 3)
    $cell.Value = "CheckBestPractices\" + $cell.Value2
}

$wsWorkflow.Range("E4").Select()

# --- Project sheet ----------------------------------------------------
$wsProject = $wb.Worksheets.Item("Project")

$cell = $wsProject.Cells.Item(2, 3)
$cell.Value = "CheckBestPractices\" + $cell.Value2

$wsProject.Range("E1").Value = "Action"
$wsProject.Range("E2").Value = "Double check"

$wsProject.Range("E3").Select()
